$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text so numeric-looking strings (e.g. "248.92")
# are preserved verbatim instead of being parsed into floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.209.44"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.064.13"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "248.92"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "0.667"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "57.33"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "16.31"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "0.924"
$ws.Range("E13").Value = "  +14.39%  "
$ws.Range("D14").Value = "2.365.17"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "5.78"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "2.067.58"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "18.76"
$ws.Range("E17").Value = "  +11.13%  "
$ws.Range("D18").Value = "37.207.37"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "75.03"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "5.50"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "238.18"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +4.52%  "
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").Value = "170.05"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "20.29"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "5.17"
$ws.Range("E30").Value = "  +8.80%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D33").Value = "4.67"
$ws.Range("E33").Value = "  +4.25%  "
$ws.Range("D34").Value = "0.0894"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").Value = "1.78"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  +14.91%  "
$ws.Range("D40").Value = "3.11"
$ws.Range("E40").Value = "  +8.12%  "
$ws.Range("E41").Value = "  -10.91%  "
$ws.Range("D42").Value = "17.68"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "1.17"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "96.89"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "1.278.00"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "2.87"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "2.252.14"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "44.55"
$ws.Range("E51").Value = "  +1.32%  "

# Restore the original (default) style on the Price column now that the
# values are locked in as text, so no stray number-format style lingers.
$priceRange.Style = "Normal"
